$wb = $excel.ActiveWorkbook

# Sheet: "accuracy period3testperiod3"  (sheet2.xml)
$ws = $wb.Worksheets.Item("accuracy period3testperiod3")
$ws.Range("A1").Value = 0.23255813953488369
$ws.Range("D1").Value = 0.55348837209302326
$ws.Range("A2").Value = 0.76744186046511631
$ws.Range("A3").Value = 0.23255813953488369
$ws.Range("A4").Value = 0.76744186046511631
$ws.Range("A5").Value = 0.76744186046511631

# Sheet: "accuracy period3testperiod4"  (sheet3.xml)
$ws = $wb.Worksheets.Item("accuracy period3testperiod4")
$ws.Range("A1").Value = 0.76744186046511631
$ws.Range("D1").Value = 0.76744186046511631
$ws.Range("A2").Value = 0.76744186046511631
$ws.Range("A3").Value = 0.76744186046511631
$ws.Range("A4").Value = 0.76744186046511631
$ws.Range("A5").Value = 0.76744186046511631

# Sheet: "accuracy period3testperiod5"  (sheet4.xml)
$ws = $wb.Worksheets.Item("accuracy period3testperiod5")
$ws.Range("A1").Value = 0.76744186046511631
$ws.Range("D1").Value = 0.76744186046511631
$ws.Range("A2").Value = 0.76744186046511631
$ws.Range("A3").Value = 0.76744186046511631
$ws.Range("A4").Value = 0.76744186046511631
$ws.Range("A5").Value = 0.76744186046511631

# Restore the active sheet / selection as in the original workbook view
$wb.Worksheets.Item("accuracy period3testperiod5").Activate()
